# Reproduce the author's edit: on the "DELIVER" sheet, rows 1-4 (which were
# blank placeholder rows above the actual data that started at row 5) are
# deleted, shifting everything up by 4 rows. Tables, merged cells and data
# validation ranges all move with the rows automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DELIVER")
$ws.Activate()

# Delete the first 4 (empty) rows - this shifts rows 5.. up to 1..
$ws.Rows("1:4").Delete()

# The final row selected by the author after editing.
$ws.Range("B10").Select()
